# edit.ps1
# Applies the "fix nan in required data" changes:
#  1. "Collect Redress Kit" sheet: update Required (C) / Need to order (G) / Reserved (H)
#     columns for existing rows, and append new rows (14-34) for additional redress kits.
#  2. "Store" sheet: update on-hand Quantity (B) for several items.

$wb = $excel.ActiveWorkbook
$wsKit = $wb.Worksheets.Item("Collect Redress Kit")
$wsStore = $wb.Worksheets.Item("Store")

$kitRows = @(
    @{Row=2; A='REDRESS0950'; B=12; C=1; D=0; E='L043'; F='Ball bearing'; G=0; H=2},
    @{Row=3; A='REDRESS0950'; B=12; C=1; D=0; E='L107'; F='Needle bearing, NK10/16'; G=0; H=3},
    @{Row=4; A='REDRESS0950'; B=12; C=1; D=0; E='T009'; F='O-ring 18,77 x 1,78 - V70'; G=0; H=1},
    @{Row=5; A='REDRESS0950'; B=12; C=1; D=0; E='T006'; F='O-ring 37,82 x 1,78 - V70'; G=0; H=1},
    @{Row=6; A='REDRESS0950'; B=12; C=1; D=0; E='T349'; F='O-ring 6  x 1,0 - V70'; G=0; H=2},
    @{Row=7; A='REDRESS0950'; B=12; C=1; D=0; E='HT25046'; F='Seal'; G=0; H=48},
    @{Row=8; A='REDRESS0950'; B=12; C=1; D=0; E='T354'; F='O-ring 5,2 x 0,6 - FPM75'; G=34; H=14},
    @{Row=9; A='REDRESS0950'; B=12; C=1; D=0; E='T209'; F='O-ring 5,50 x 1,00 - V70'; G=10; H=38},
    @{Row=10; A='REDRESS0950'; B=12; C=1; D=0; E='F585'; F='Compression spring SIF-12337'; G=0; H=12},
    @{Row=11; A='REDRESS0950'; B=12; C=1; D=0; E='F211'; F='Compression spring 20750'; G=24; H=0},
    @{Row=12; A='REDRESS0950'; B=12; C=1; D=0; E='HT24035'; F='Retaining ring f HS24010 valve'; G=20; H=4},
    @{Row=13; A='REDRESS0950'; B=12; C=1; D=0; E='R081'; F='Retaining ring A6'; G=0; H=2},
    @{Row=14; A='REDRESS0764'; B=96; C=0; D=0; E='R053'; F='Retaining ring INCONEL X-750'; G=0; H=2},
    @{Row=15; A='REDRESS0764'; B=96; C=0; D=0; E='R130'; F='Retaining ring 218B'; G=0; H=2},
    @{Row=16; A='REDRESS0764'; B=96; C=0; D=0; E='T135'; F='O-ring 48,00 x 1,20 - V80'; G=0; H=4},
    @{Row=17; A='REDRESS0764'; B=96; C=0; D=0; E='RT21014'; F='Screw M2,5x4 mm kval.12,9'; G=0; H=20},
    @{Row=18; A='REDRESS0764'; B=96; C=0; D=0; E='T077'; F='PTFE-Ring ø3,7/2,7x0,9'; G=0; H=32},
    @{Row=19; A='REDRESS0764'; B=96; C=0; D=0; E='T081'; F='O-ring 30,0 x 1,0 - V70'; G=0; H=8},
    @{Row=20; A='REDRESS0764'; B=96; C=0; D=0; E='T006'; F='O-ring 37,82 x 1,78 - V80'; G=0; H=4},
    @{Row=21; A='REDRESS0764'; B=96; C=0; D=0; E='R023'; F='Retaining ring'; G=0; H=4},
    @{Row=22; A='REDRESS0764'; B=96; C=0; D=0; E='T095'; F='O-ring 3,00 x 1,00 - V70'; G=0; H=8},
    @{Row=23; A='REDRESS0764'; B=96; C=0; D=0; E='T020'; F='O-ring 6,07 x 1,78 - V75'; G=0; H=4},
    @{Row=24; A='REDRESS0764'; B=96; C=0; D=0; E='T049'; F='O-ring 14,00 x 1,78 - V75'; G=0; H=8},
    @{Row=25; A='REDRESS0180'; B=74; C=1; D=1; E='T198'; F='O-ring ø1,78 - V70 L=1423'; G=0; H=1},
    @{Row=26; A='REDRESS0180'; B=74; C=1; D=1; E='T001'; F='O-ring 7,66 x 1,78 - V70'; G=0; H=1},
    @{Row=27; A='REDRESS0180'; B=74; C=1; D=1; E='T070'; F='O-ring 4,50 x 1,00 - V70'; G=0; H=12},
    @{Row=28; A='REDRESS0180'; B=74; C=1; D=1; E='T021'; F='O-ring 10,82 x 1,78 - V70'; G=0; H=4},
    @{Row=29; A='REDRESS0180'; B=74; C=1; D=1; E='T071'; F='Glide ring'; G=0; H=4},
    @{Row=30; A='REDRESS0180'; B=74; C=1; D=1; E='KT22024'; F='Screw M4'; G=0; H=10},
    @{Row=31; A='REDRESS0180'; B=74; C=1; D=1; E='T006'; F='O-ring 37,82 x 1,78 - V70'; G=0; H=1},
    @{Row=32; A='REDRESS0180'; B=74; C=1; D=1; E='T015'; F='O-ring 5,28 x 1,78-V70'; G=0; H=4},
    @{Row=33; A='REDRESS0180'; B=74; C=1; D=1; E='T025'; F='O-ring 9,25 x 1,78 - V70'; G=0; H=6},
    @{Row=34; A='REDRESS0180'; B=74; C=1; D=1; E='T020'; F='O-ring 6,07 x 1,78 - V70'; G=0; H=3}
)

$storeRows = @(
    @{Row=787; B=8},
    @{Row=846; B=26},
    @{Row=909; B=9},
    @{Row=935; B=7},
    @{Row=1154; B=0},
    @{Row=1165; B=24},
    @{Row=1169; B=5},
    @{Row=1185; B=26},
    @{Row=1515; B=32},
    @{Row=1773; B=45},
    @{Row=1782; B=21},
    @{Row=1790; B=11},
    @{Row=1801; B=70},
    @{Row=1811; B=44},
    @{Row=1814; B=72},
    @{Row=1819; B=35},
    @{Row=1847; B=13},
    @{Row=1870; B=0},
    @{Row=1871; B=16},
    @{Row=1875; B=323},
    @{Row=1881; B=45},
    @{Row=1896; B=1},
    @{Row=1915; B=16},
    @{Row=1951; B=3},
    @{Row=2025; B=7}
)

# --- Update / create rows on "Collect Redress Kit" ---

# Make sure new rows (14-34) inherit the same formatting as the existing data rows
# (style used by row 13) before writing values into them.
$wsKit.Range("A13:H13").Copy() | Out-Null
$wsKit.Range("A14:H34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($r in $kitRows) {
    $row = $r.Row
    $wsKit.Cells.Item($row, 1).Value2 = $r.A
    $wsKit.Cells.Item($row, 2).Value2 = $r.B
    $wsKit.Cells.Item($row, 3).Value2 = $r.C
    $wsKit.Cells.Item($row, 4).Value2 = $r.D
    $wsKit.Cells.Item($row, 5).Value2 = $r.E
    $wsKit.Cells.Item($row, 6).Value2 = $r.F
    $wsKit.Cells.Item($row, 7).Value2 = $r.G
    $wsKit.Cells.Item($row, 8).Value2 = $r.H
}

# --- Update on-hand quantities on "Store" ---

foreach ($r in $storeRows) {
    $wsStore.Cells.Item($r.Row, 2).Value2 = $r.B
}

Write-Output "Updated $($kitRows.Count) rows on 'Collect Redress Kit' and $($storeRows.Count) rows on 'Store'."
